# Update the Reaction_number (column C) values on both the NBR and BAR
# sheets to reflect the newly re-run sensitivity analysis (new ecoli
# studies / NB ratio analysis).

$wb = $excel.ActiveWorkbook

$wsNBR = $wb.Worksheets.Item("NBR")
$wsBAR = $wb.Worksheets.Item("BAR")

# New column C values for rows 2..20 (Cutoff 0..18 / Reaction_number 1..19)
$nbrValues = @(526, 525, 512, 510, 504, 498, 491, 480, 478, 472, 472, 469, 461, 456, 0, 441, 445, 440, 436)
$barValues = @(645, 648, 646, 643, 644, 641, 642, 637, 637, 634, 633, 636, 633, 633, 0, 626, 627, 621, 620)

for ($i = 0; $i -lt $nbrValues.Length; $i++) {
    $row = $i + 2
    $wsNBR.Cells.Item($row, 3).Value = $nbrValues[$i]
    $wsBAR.Cells.Item($row, 3).Value = $barValues[$i]
}
